$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a copy of row 16 at position 17 (shifts existing rows 17-44 down to 18-45)
# so the new row inherits the same formatting (e.g. date style on column D).
$ws.Rows.Item(16).Copy()
$ws.Rows.Item(17).Insert()

$newRow = 17

$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 44483
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 100112052
$ws.Cells.Item($newRow, 7).Value = "Albahaca"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 700
$ws.Cells.Item($newRow, 11).Value = 3500
$ws.Cells.Item($newRow, 12).Value = 4000
$ws.Cells.Item($newRow, 13).Value = 3750
$ws.Cells.Item($newRow, 14).Value = "`$/paquete"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 3750
$ws.Cells.Item($newRow, 17).Value = 1
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
